$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename to machine-friendly column names ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Title-case the "de/del/el/la/los/y" connector words in municipality/state names ---
$ws.Range('B6').Value = 'Rincón De Romos'
$ws.Range('B7').Value = 'San José De Gracia'
$ws.Range('B23').Value = 'Amatenango De La Frontera'
$ws.Range('B32').Value = 'Comitán De Domínguez'
$ws.Range('B48').Value = 'Montecristo De Guerrero'
$ws.Range('B51').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B55').Value = 'San Cristóbal De Las Casas'
$ws.Range('B77').Value = 'Hidalgo Del Parral'
$ws.Range('B100').Value = 'San Juan De Sabinas'
$ws.Range('A107').Value = 'Ciudad De México'
$ws.Range('B111').Value = 'Cuajimalpa De Morelos'
$ws.Range('B136').Value = 'San Juan Del Río'
$ws.Range('B559').Value = 'San Juan Del Río'
$ws.Range('A141').Value = 'Estado De México'
$ws.Range('B141').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B143').Value = 'Almoloya De Alquisiras'
$ws.Range('B152').Value = 'Coacalco De Berriozábal'
$ws.Range('B154').Value = 'Ecatepec De Morelos'
$ws.Range('B158').Value = 'Naucalpan De Juárez'
$ws.Range('B160').Value = 'San Felipe Del Progreso'
$ws.Range('B166').Value = 'Tenango Del Aire'
$ws.Range('B170').Value = 'Tlalnepantla De Baz'
$ws.Range('B174').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B182').Value = 'Apaseo El Alto'
$ws.Range('B183').Value = 'Apaseo El Grande'
$ws.Range('B189').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B202').Value = 'San Diego De La Unión'
$ws.Range('B204').Value = 'San Francisco Del Rincón'
$ws.Range('B206').Value = 'San Luis De La Paz'
$ws.Range('B207').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B208').Value = 'Silao De La Victoria'
$ws.Range('B212').Value = 'Valle De Santiago'
$ws.Range('B218').Value = 'Acapulco De Juárez'
$ws.Range('B222').Value = 'Atoyac De Álvarez'
$ws.Range('B223').Value = 'Ayutla De Los Libres'
$ws.Range('B225').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B227').Value = 'Coyuca De Benítez'
$ws.Range('B228').Value = 'Coyuca De Catalán'
$ws.Range('B231').Value = 'Cuetzala Del Progreso'
$ws.Range('B232').Value = 'Cutzamala De Pinzón'
$ws.Range('B236').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B237').Value = 'Iguala De La Independencia'
$ws.Range('B240').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B250').Value = 'Taxco De Alarcón'
$ws.Range('B251').Value = 'Técpan De Galeana'
$ws.Range('B256').Value = 'Tlapa De Comonfort'
$ws.Range('B261').Value = 'Agua Blanca De Iturbide'
$ws.Range('B265').Value = 'Atotonilco De Tula'
$ws.Range('B266').Value = 'Atotonilco El Grande'
$ws.Range('B270').Value = 'Cuautepec De Hinojosa'
$ws.Range('B272').Value = 'Huasca De Ocampo'
$ws.Range('B273').Value = 'Huejutla De Reyes'
$ws.Range('B281').Value = 'Mixquiahuala De Juárez'
$ws.Range('B282').Value = 'Pachuca De Soto'
$ws.Range('B288').Value = 'Tenango De Doria'
$ws.Range('B289').Value = 'Tepehuacán De Guerrero'
$ws.Range('B290').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B291').Value = 'Tezontepec De Aldama'
$ws.Range('B296').Value = 'Tula De Allende'
$ws.Range('B297').Value = 'Tulancingo De Bravo'
$ws.Range('B299').Value = 'Zacualtipán De Ángeles'
$ws.Range('B305').Value = 'Atotonilco El Alto'
$ws.Range('B306').Value = 'Autlán De Navarro'
$ws.Range('B311').Value = 'Encarnación De Díaz'
$ws.Range('B317').Value = 'Lagos De Moreno'
$ws.Range('B322').Value = 'Ojuelos De Jalisco'
$ws.Range('B323').Value = 'Talpa De Allende'
$ws.Range('B324').Value = 'Tamazula De Gordiano'
$ws.Range('B326').Value = 'Tepatitlán De Morelos'
$ws.Range('B327').Value = 'Valle De Juárez'
$ws.Range('B329').Value = 'Yahualica De González Gallo'
$ws.Range('B339').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B372').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B391').Value = 'Coatlán Del Río'
$ws.Range('B398').Value = 'Puente De Ixtla'
$ws.Range('B400').Value = 'Tlaltizapán De Zapata'
$ws.Range('B405').Value = 'Ixtlán Del Río'
$ws.Range('B406').Value = 'Santa María Del Oro'
$ws.Range('B420').Value = 'Mier Y Noriega'
$ws.Range('B421').Value = 'Montemorelos'
$ws.Range('B423').Value = 'San Nicolás De Los Garza'
$ws.Range('B428').Value = 'Guevea De Humboldt'
$ws.Range('B429').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B430').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B431').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B436').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B437').Value = 'Oaxaca De Juárez'
$ws.Range('B438').Value = 'Ocotlán De Morelos'
$ws.Range('B439').Value = 'Putla Villa De Guerrero'
$ws.Range('B446').Value = 'San José Del Progreso'
$ws.Range('B448').Value = 'San Juan Bautista Lo De Soto'
$ws.Range('B457').Value = 'San Miguel Del Puerto'
$ws.Range('B469').Value = 'Santa María Jalapa Del Marqués'
$ws.Range('B478').Value = 'Santo Domingo De Morelos'
$ws.Range('B483').Value = 'Tataltepec De Valdés'
$ws.Range('B484').Value = 'Teotitlán De Flores Magón'
$ws.Range('B485').Value = 'Tezoatlán De Segura Y Luna'
$ws.Range('B486').Value = 'Tlacolula De Matamoros'
$ws.Range('B487').Value = 'Villa De Etla'
$ws.Range('B488').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B489').Value = 'Villa Sola De Vega'
$ws.Range('B490').Value = 'Zimatlán De Álvarez'
$ws.Range('B502').Value = 'Cuayuca De Andrade'
$ws.Range('B503').Value = 'Cuetzalan Del Progreso'
$ws.Range('B510').Value = 'Izúcar De Matamoros'
$ws.Range('B517').Value = 'Palmar De Bravo'
$ws.Range('B524').Value = 'San Salvador El Seco'
$ws.Range('B531').Value = 'Tepexi De Rodríguez'
$ws.Range('B532').Value = 'Tetela De Ocampo'
$ws.Range('B533').Value = 'Teteles De Avila Castillo'
$ws.Range('B536').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B544').Value = 'Xochitlán De Vicente Suárez'
$ws.Range('B549').Value = 'Amealco De Bonfil'
$ws.Range('B554').Value = 'Jalpan De Serra'
$ws.Range('B555').Value = 'Landa De Matamoros'
$ws.Range('B557').Value = 'Pinal De Amoles'
$ws.Range('B568').Value = 'Armadillo De Los Infante'
$ws.Range('B569').Value = 'Axtla De Terrazas'
$ws.Range('B574').Value = 'Ciudad Del Maíz'
$ws.Range('B583').Value = 'Mexquitic De Carmona'
$ws.Range('B588').Value = 'San Ciro De Acosta'
$ws.Range('B594').Value = 'Santa María Del Río'
$ws.Range('B599').Value = 'Tanquián De Escobedo'
$ws.Range('B603').Value = 'Villa De Arista'
$ws.Range('B604').Value = 'Villa De Arriaga'
$ws.Range('B605').Value = 'Villa De Guadalupe'
$ws.Range('B606').Value = 'Villa De Ramos'
$ws.Range('B607').Value = 'Villa De Reyes'
$ws.Range('B640').Value = 'Jalpa De Méndez'
$ws.Range('B661').Value = 'Soto La Marina'
$ws.Range('B670').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B671').Value = 'Sanctórum De Lázaro Cárdenas'
$ws.Range('B684').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B687').Value = 'Amatlán De Los Reyes'
$ws.Range('B695').Value = 'Cazones De Herrera'
$ws.Range('B705').Value = 'Cosamaloapan De Carpio'
$ws.Range('B706').Value = 'Cosautlán De Carvajal'
$ws.Range('B715').Value = 'Hueyapan De Ocampo'
$ws.Range('B716').Value = 'Ignacio De La Llave'
$ws.Range('B718').Value = 'Ixhuatlán Del Café'
$ws.Range('B723').Value = 'Juchique De Ferrer'
$ws.Range('B726').Value = 'Lerdo De Tejada'
$ws.Range('B728').Value = 'Martínez De La Torre'
$ws.Range('B730').Value = 'Medellín De Bravo'
$ws.Range('B738').Value = 'Paso De Ovejas'
$ws.Range('B739').Value = 'Paso Del Macho'
$ws.Range('B742').Value = 'Poza Rica De Hidalgo'
$ws.Range('B748').Value = 'Sayula De Alemán'
$ws.Range('B749').Value = 'Soledad De Doblado'
$ws.Range('B752').Value = 'Tatahuicapan De Juárez'
$ws.Range('B784').Value = 'Concepción Del Oro'
$ws.Range('B796').Value = 'Nochistlán De Mejía'
$ws.Range('B803').Value = 'Tlaltenango De Sánchez Román'

# --- Floating point recalculation artifact on D776 ---
$ws.Range('D776').Value = 0.09769094138543516

# --- Remove trailing footer/metadata rows 811-815 (row 810 was already blank) ---
$ws.Rows('811:815').Delete()
